# katalog.xlsx — "Add files via upload"
#
# Adds a discount ("Rabatt in %" / p_rabatt) input row to the
# "Brix_Gel_Stab" worksheet and folds that discount factor into the
# sheet's price formula documentation string, mirroring the pattern that
# already exists on "Eigen_Edelstahl".

$wb = $excel.ActiveWorkbook

# --- Incidental cursor-position change on another sheet picked up by the
# same autosave. Do this first so the final active tab below ends up being
# "Brix_Gel_Stab" again (matches tabSelected/activeTab in the target). ---
$ws2 = $wb.Worksheets.Item("Brix_Zaun_Stab")
$ws2.Range("E9").Select() | Out-Null

# --- Brix_Gel_Stab: insert a new row 9 (Zahl / Rabatt in % / p_rabatt / 15)
# pushing the old totals row (Preis/Gesamtpreis/Endpreis) down to row 10. ---
$ws = $wb.Worksheets.Item("Brix_Gel_Stab")
$ws.Activate() | Out-Null

$ws.Rows.Item(9).Insert() | Out-Null

$ws.Range("A9").Value = "Zahl"
$ws.Range("B9").Value = "Rabatt in %"
$ws.Range("C9").Value = "p_rabatt"
$ws.Range("D9").Value = 15

# Update the "Endpreis" formula documentation cell (now on row 10) so the
# discount factor is applied.
$ws.Range("E10").Value = "((P_Modell * L * F_Faktor * P_Handlauf) + ((math.ceil(L/1.3)+1) * P_Steher * F_Faktor) + (Ecken * 95) + (L * P_Arbeit) * ( 1 - (p_rabatt / 100))"

# Leave the cursor where the author last left it on this sheet.
$ws.Range("D15").Select() | Out-Null
